$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 383.33334
$ws.Range("J32").Value = 375
$ws.Range("L32").Value = 375
$ws.Range("N32").Value = -1027
# Row 43
$ws.Range("H43").Value = 1918.04
$ws.Range("I43").Value = 1790.8667
$ws.Range("K43").Value = 1790.8667
$ws.Range("M43").Value = -1721.8667
# Row 132
$ws.Range("H132").Value = 5819838.5
$ws.Range("I132").Value = 6950529
$ws.Range("J132").Value = 4858.7144
$ws.Range("K132").Value = 20851587
$ws.Range("L132").Value = 14576.1432
$ws.Range("M132").Value = -20849057
$ws.Range("N132").Value = -19636.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 2605.2327
$ws.Range("I132").Value = 2694.275
$ws.Range("J132").Value = 1418
$ws.Range("K132").Value = 8082.825000000001
$ws.Range("L132").Value = 4254
$ws.Range("M132").Value = -5552.825000000001
$ws.Range("N132").Value = -9314

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2927.1143
$ws.Range("J134").Value = 2600
$ws.Range("L134").Value = 7800
$ws.Range("N134").Value = -12870

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1356.7778
$ws.Range("I16").Value = 1333.3334
$ws.Range("J16").Value = 1368.5
$ws.Range("K16").Value = 1333.3334
$ws.Range("L16").Value = 1368.5
$ws.Range("M16").Value = -1046.3334
$ws.Range("N16").Value = -1942.5
# Row 31
$ws.Range("H31").Value = 21877.64
$ws.Range("I31").Value = 1109.0883
$ws.Range("J31").Value = 37926.066
$ws.Range("K31").Value = 1109.0883
$ws.Range("L31").Value = 37926.066
$ws.Range("M31").Value = -814.0882999999999
$ws.Range("N31").Value = -38516.066
# Row 34
$ws.Range("H34").Value = 21877.64
$ws.Range("I34").Value = 1109.0883
$ws.Range("J34").Value = 37926.066
$ws.Range("K34").Value = 1109.0883
$ws.Range("L34").Value = 37926.066
$ws.Range("M34").Value = -907.0882999999999
$ws.Range("N34").Value = -38330.066
# Row 36
$ws.Range("H36").Value = 9000
$ws.Range("I36").Value = 9000
$ws.Range("K36").Value = 9000
$ws.Range("M36").Value = -8612
# Row 40
$ws.Range("H40").Value = 9000
$ws.Range("I40").Value = 9000
$ws.Range("K40").Value = 9000
$ws.Range("M40").Value = -8840
# Row 62
$ws.Range("H62").Value = 2091
$ws.Range("I62").Value = 1038.3334
$ws.Range("J62").Value = 2542.1428
$ws.Range("K62").Value = 1038.3334
$ws.Range("L62").Value = 2542.1428
$ws.Range("M62").Value = -414.3334
$ws.Range("N62").Value = -3790.1428
# Row 65
$ws.Range("H65").Value = 2091
$ws.Range("I65").Value = 1038.3334
$ws.Range("J65").Value = 2542.1428
$ws.Range("K65").Value = 5191.666999999999
$ws.Range("L65").Value = 12710.714
$ws.Range("M65").Value = -2071.666999999999
$ws.Range("N65").Value = -18950.714
# Row 94
$ws.Range("H94").Value = 1195
$ws.Range("I94").Value = 1056
$ws.Range("J94").Value = 1229.75
$ws.Range("K94").Value = 1056
$ws.Range("L94").Value = 1229.75
$ws.Range("M94").Value = -605
$ws.Range("N94").Value = -2131.75
# Row 113
$ws.Range("H113").Value = 1356.7778
$ws.Range("I113").Value = 1333.3334
$ws.Range("J113").Value = 1368.5
$ws.Range("K113").Value = 1333.3334
$ws.Range("L113").Value = 1368.5
$ws.Range("M113").Value = 836.6666
$ws.Range("N113").Value = -5708.5
# Row 132
$ws.Range("H132").Value = 68187020
$ws.Range("I132").Value = 71434690
$ws.Range("K132").Value = 214304070
$ws.Range("M132").Value = -214301540

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 171.66667
$ws.Range("I2").Value = 239
$ws.Range("K2").Value = 1434
$ws.Range("M2").Value = -1321
# Row 18
$ws.Range("H18").Value = 357.6842
$ws.Range("I18").Value = 327.6875
$ws.Range("K18").Value = 983.0625
$ws.Range("M18").Value = -814.0625
# Row 118
$ws.Range("H118").Value = 1746.6666
$ws.Range("I118").Value = 1096
$ws.Range("K118").Value = 3288
$ws.Range("M118").Value = -2045
# Row 131
$ws.Range("H131").Value = 807.97
$ws.Range("J131").Value = 865.6477
$ws.Range("L131").Value = 2596.9431
$ws.Range("N131").Value = -12676.9431

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 100104750
$ws.Range("J80").Value = 1900
$ws.Range("L80").Value = 1900
$ws.Range("N80").Value = -3896
# Row 83
$ws.Range("H83").Value = 100104750
$ws.Range("J83").Value = 1900
$ws.Range("L83").Value = 9500
$ws.Range("N83").Value = -19484
# Row 132
$ws.Range("H132").Value = 2370.0588
$ws.Range("I132").Value = 1708.2727
$ws.Range("J132").Value = 3583.3333
$ws.Range("K132").Value = 5124.8181
$ws.Range("L132").Value = 10749.9999
$ws.Range("M132").Value = -2594.8181
$ws.Range("N132").Value = -15809.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3397.45
$ws.Range("I7").Value = 1740.2
$ws.Range("J7").Value = 5054.7
$ws.Range("K7").Value = 1740.2
$ws.Range("L7").Value = 5054.7
$ws.Range("M7").Value = -1628.2
$ws.Range("N7").Value = -5278.7
# Row 22
$ws.Range("H22").Value = 839.6667
$ws.Range("I22").Value = 1124.75
$ws.Range("J22").Value = 772.58826
$ws.Range("K22").Value = 1124.75
$ws.Range("L22").Value = 772.58826
$ws.Range("M22").Value = -829.75
$ws.Range("N22").Value = -1362.58826
# Row 27
$ws.Range("H27").Value = 839.6667
$ws.Range("I27").Value = 1124.75
$ws.Range("J27").Value = 772.58826
$ws.Range("K27").Value = 1124.75
$ws.Range("L27").Value = 772.58826
$ws.Range("M27").Value = -1017.75
$ws.Range("N27").Value = -986.58826
# Row 55
$ws.Range("H55").Value = 447.34784
$ws.Range("J55").Value = 539.8570999999999
$ws.Range("L55").Value = 539.8570999999999
$ws.Range("N55").Value = -885.8570999999999
# Row 82
$ws.Range("H82").Value = 1167.3846
$ws.Range("I82").Value = 908
$ws.Range("J82").Value = 1329.5
$ws.Range("K82").Value = 908
$ws.Range("L82").Value = 1329.5
$ws.Range("M82").Value = -547
$ws.Range("N82").Value = -2051.5
# Row 85
$ws.Range("H85").Value = 1167.3846
$ws.Range("I85").Value = 908
$ws.Range("J85").Value = 1329.5
$ws.Range("K85").Value = 908
$ws.Range("L85").Value = 1329.5
$ws.Range("M85").Value = 340
$ws.Range("N85").Value = -3825.5
# Row 126
$ws.Range("H126").Value = 3397.45
$ws.Range("I126").Value = 1740.2
$ws.Range("J126").Value = 5054.7
$ws.Range("K126").Value = 5220.6
$ws.Range("L126").Value = 15164.1
$ws.Range("M126").Value = -2750.6
$ws.Range("N126").Value = -20104.1
# Row 132
$ws.Range("H132").Value = 3160.3635
$ws.Range("I132").Value = 3466.4075
$ws.Range("J132").Value = 1783.1666
$ws.Range("K132").Value = 10399.2225
$ws.Range("L132").Value = 5349.4998
$ws.Range("M132").Value = -7869.2225
$ws.Range("N132").Value = -10409.4998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 7000.1816
$ws.Range("I132").Value = 10184.4
$ws.Range("J132").Value = 4346.6665
$ws.Range("K132").Value = 30553.2
$ws.Range("L132").Value = 13039.9995
$ws.Range("M132").Value = -28023.2
$ws.Range("N132").Value = -18099.9995
